$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write a value into a cell as TEXT, even if it "looks like" a
# number (e.g. "1.001"), without leaving any NumberFormat/style residue on
# the destination cell.
#
# Plain `$range.Value = "1.001"` lets Excel's usual type-inference kick in
# and the cell is silently stored as the *number* 1.001 (losing the
# trailing-zero/precision formatting of the original text). Pre-stamping
# the destination with NumberFormat "@" avoids that but permanently tags
# the cell with a new style index, which the source workbook never had.
#
# Instead we stamp a scratch cell (Z1, well outside the used range) as
# Text, put the value there, and Copy/PasteSpecial *values only* into the
# real destination - that carries over the literal text without copying
# the "@" number format, then we clear the scratch cell again.
$helper = $ws.Range("Z1")
function Set-TextValue {
    param($range, [string]$value)
    $helper.NumberFormat = "@"
    $helper.Value = $value
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
}

# Row => updated Coin(B) / Link(C) / Price(D) / Volume-1h(E) values.
# B/C are only present where the coin in that row changed (TrustWalletToken
# <-> Frax swapped rows 36/37, TheSandbox <-> Algorand swapped rows 41/42).
$rows = @(
    @{Row=2; D='28.234.73'; E='  +2.36%  '},
    @{Row=3; D='1.878.79'; E='  +1.55%  '},
    @{Row=4; D='1.001'; E='  -0.20%  '},
    @{Row=5; D='315.57'; E='  +0.44%  '},
    @{Row=6; D='1.001'; E='  -0.20%  '},
    @{Row=7; D='0.4312'; E='  +1.50%  '},
    @{Row=8; D='0.3706'; E='  +1.89%  '},
    @{Row=9; D='0.07419'; E='  +1.66%  '},
    @{Row=10; D='0.8837'; E='  +1.19%  '},
    @{Row=11; D='21.17'; E='  +2.22%  '},
    @{Row=12; D='1.894.15'; E='  -0.71%  '},
    @{Row=13; D='5.497'; E='  +2.84%  '},
    @{Row=14; D='6.642'; E='  +1.82%  '},
    @{Row=15; D='0.06987'; E='  +0.94%  '},
    @{Row=16; D='1.003'; E='  -0.04%  '},
    @{Row=17; D='81.31'; E='  +3.00%  '},
    @{Row=18; D='0.000009154'; E='  +3.16%  '},
    @{Row=19; D='1.001'; E='  -0.25%  '},
    @{Row=20; D='15.63'; E='  +1.47%  '},
    @{Row=21; D='28.277.34'; E='  +2.40%  '},
    @{Row=22; D='5.093'; E='  +1.69%  '},
    @{Row=23; D='10.96'; E='  +3.25%  '},
    @{Row=24; D='2.124.83'; E='  +0.15%  '},
    @{Row=25; D='1.978'; E='  -0.32%  '},
    @{Row=26; D='154.19'; E='  +0.36%  '},
    @{Row=27; D='18.73'; E='  -1.48%  '},
    @{Row=28; D='5.434'; E='  +3.14%  '},
    @{Row=29; D='117.80'; E='  -3.03%  '},
    @{Row=30; D='1.874'; E='  -1.71%  '},
    @{Row=31; D='0.08968'; E='  +0.58%  '},
    @{Row=32; D='0.7955'; E='  +4.44%  '},
    @{Row=33; D='4.730'; E='  +3.44%  '},
    @{Row=34; D='1.185'; E='  +7.72%  '},
    @{Row=35; D='2.954'; E='  +0.49%  '},
    @{Row=36; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.131'; E='  +3.57%  '},
    @{Row=37; B='Frax'; C='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D='1.001'; E='  -0.19%  '},
    @{Row=38; D='0.05462'; E='  +1.84%  '},
    @{Row=39; D='0.01971'; E='  +1.53%  '},
    @{Row=40; D='2.889'; E='  +3.03%  '},
    @{Row=41; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5177'; E='  +1.38%  '},
    @{Row=42; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1694'; E='  +2.70%  '},
    @{Row=43; D='6.882'; E='  -0.89%  '},
    @{Row=44; D='8.634'; E='  +4.14%  '},
    @{Row=45; D='10.58'; E='  +1.36%  '},
    @{Row=46; D='0.06585'; E='  +0.40%  '},
    @{Row=47; D='0.4767'; E='  +0.37%  '},
    @{Row=48; D='106.05'; E='  +1.47%  '},
    @{Row=49; D='1.000'; E='  -0.24%  '},
    @{Row=50; D='1.658'; E='  +2.06%  '},
    @{Row=51; D='1.850'; E='  +5.90%  '}
)

foreach ($r in $rows) {
    if ($r.ContainsKey("B")) { $ws.Cells.Item($r.Row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($r.Row, 5).Value = $r.E }
}
